$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 used to be the first "Contact" row; turn it into "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely
$ws1.Range("A11").EntireRow.Delete()
